# Update the cryptos list with freshly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    # Decide whether Excel would reinterpret this text as a number (e.g.
    # "592.23"). If so, force the cell format to Text first so the stored
    # value stays a string, matching the sheet's original inline-string
    # cell type. Values that already aren't numeric-looking (contain
    # extra separators, "%", letters, subscripts, etc.) are left alone.
    $looksNumeric = $text -match '^[+-]?[0-9]+(\.[0-9]+)?$'
    if ($looksNumeric) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $text
}

# Each entry: row number, D-column price text (or $null to leave unchanged),
# E-column volume/percentage text (always "  +x.xx%  " padded with two spaces).
$updates = @(
    @{ Row = 2;  D = "68.293.91";   E = "  +0.80%  " }
    @{ Row = 3;  D = "2.542.62";    E = "  +0.58%  " }
    @{ Row = 4;  D = $null;         E = "  +0.03%  " }
    @{ Row = 5;  D = "592.23";      E = "  -0.07%  " }
    @{ Row = 6;  D = "175.13";      E = "  -0.84%  " }
    @{ Row = 7;  D = $null;         E = "  +0.00%  " }
    @{ Row = 8;  D = $null;         E = "  -0.92%  " }
    @{ Row = 9;  D = "2.542.04";    E = "  +0.65%  " }
    @{ Row = 10; D = "0.137";       E = "  -2.53%  " }
    @{ Row = 11; D = "0.166";       E = "  +1.76%  " }
    @{ Row = 12; D = $null;         E = "  +0.28%  " }
    @{ Row = 13; D = "5.02";        E = "  -2.61%  " }
    @{ Row = 14; D = "26.56";       E = "  -1.04%  " }
    @{ Row = 15; D = "2.989.88";    E = "  +0.10%  " }
    @{ Row = 16; D = "0.0000176";   E = "  -1.01%  " }
    @{ Row = 17; D = "68.281.39";   E = "  +0.98%  " }
    @{ Row = 18; D = "2.487.77";    E = "  -1.04%  " }
    @{ Row = 19; D = "11.93";       E = "  +4.15%  " }
    @{ Row = 20; D = $null;         E = "  -0.13%  " }
    @{ Row = 21; D = $null;         E = "  +70.68%  " }
    @{ Row = 22; D = "366.24";      E = "  +1.98%  " }
    @{ Row = 23; D = $null;         E = "  -0.75%  " }
    @{ Row = 24; D = "4.56";        E = "  -1.95%  " }
    @{ Row = 25; D = "72.07";       E = "  +1.71%  " }
    @{ Row = 26; D = $null;         E = "  -0.08%  " }
    @{ Row = 27; D = "1.89";        E = "  -4.65%  " }
    @{ Row = 28; D = "9.92";        E = "  -3.73%  " }
    @{ Row = 29; D = "2.673.83";    E = "  +0.74%  " }
    @{ Row = 30; D = "0.0₃0964";    E = "  -2.64%  " }
    @{ Row = 31; D = "535.45";      E = "  -3.55%  " }
    @{ Row = 32; D = "8.30";        E = "  +0.44%  " }
    @{ Row = 33; D = "1.30";        E = "  -3.32%  " }
    @{ Row = 34; D = $null;         E = "  +0.38%  " }
    @{ Row = 35; D = $null;         E = "  -1.26%  " }
    @{ Row = 37; D = "159.78";      E = "  +2.66%  " }
    @{ Row = 38; D = $null;         E = "  -2.40%  " }
    @{ Row = 39; D = "19.33";       E = "  +3.07%  " }
    @{ Row = 40; D = $null;         E = "  +0.24%  " }
    @{ Row = 41; D = $null;         E = "  -1.13%  " }
    @{ Row = 42; D = $null;         E = "  -0.96%  " }
    @{ Row = 43; D = $null;         E = "  -2.30%  " }
    @{ Row = 46; D = "39.46";       E = "  -1.17%  " }
    @{ Row = 47; D = "148.90";      E = "  +0.99%  " }
    @{ Row = 48; D = $null;         E = "  -0.99%  " }
    @{ Row = 49; D = $null;         E = "  -0.10%  " }
    @{ Row = 50; D = "0.0₆0277";    E = "  -0.76%  " }
    @{ Row = 51; D = $null;         E = "  +1.49%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        Set-CellValue $r 4 $u.D
    }
    Set-CellValue $r 5 $u.E
}

# Rows 44 and 45 swap their coin identity (dogwifhat <-> USDe) along with
# updated link, price and volume values.
$ws.Cells.Item(44, 2).Value = "USDe"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-CellValue 44 4 "1.00"
Set-CellValue 44 5 "  +0.43%  "

$ws.Cells.Item(45, 2).Value = "dogwifhat"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-CellValue 45 4 "2.49"
Set-CellValue 45 5 "  -1.47%  "
